# Append 5 new search-log rows (39-43) to Sheet1, per the "retested search
# engine for deployment" commit. Columns: A=prompt, C=result_1, D=result_1_info,
# E=result_2_info (B, F, G, H stay blank, matching the other logged rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 39: "what is cord blood?"
$ws.Range("A39").Value = "what is cord blood?"
$ws.Range("C39").Value = "What is cord blood? (Cord blood is blood that remains in the umbilical cord following the birth of a baby and after the umbilical cord is cut. During pregnancy, the umbilical cord functions as a lifeline between mother and child. After a baby�s delivery, the cord blood present in the umbilical cord could offer hope for the child as well as members of the family. Cord blood contains a rich source of lifesaving stem cells called Haematopoietic Stem Cells (HSCs).) (0.9998952150344849)"
$ws.Range("D39").Value = "What are cord blood stem cells? (Cord blood stem cells are also known as Haematopoietic Stem Cells (HSCs), which are responsible for replenishing blood and regenerating the immune system. These stem cells can be found within the umbilical cord following the birth of a baby and have the unique ability to differentiate into various cell types found in blood as depicted in the diagram below:     
Red blood cells - which transport oxygen
White blood cells - which produce antibodies and fight bacteria
Platelets - which assist in blood clotting                                        ) (0.7933569550514221)"
$ws.Range("E39").Value = "What can cord blood stem cells do? (Cord blood stem cells can:
Replace and regenerate damaged or diseased bone marrow.
Treat blood related cancers, solid tumours and immunodeficiency disorders.
Correct genetic defects (sibling/allogeneic transplantation).
Potentially be used for cellular therapy and regenerative medicine.erative medicine.) (0.7321342825889587)"

# Row 40: "什麼是 CDA？它是否適用於現有客戶？"
$ws.Range("A40").Value = "什麼是 CDA？它是否適用於現有客戶？"
$ws.Range("C40").Value = "儿童发展帐户CDA是什么? (CDA是儿童的一个特别储蓄账户,用于帮助积累可用于批准用途的储蓄。 儿童CDA可以在任何OCBC银行DBS银行或UOB银行分行开放。 在孩子满12岁当年的12月31日之前,您可以随时向CDA存入现金。 有关CDA的更多信息,请访问www.babybonus.msf.gov.sg, 或今天询问我们的友好顾问。 儿童如果在2008年8月17日或之后是Born,则有资格申请CDA。 新加坡公民或在12岁之前成为公民。 父母合法结婚。) (0.7435237765312195)"
$ws.Range("D40").Value = "儿童发展帐户CDA如何运作? (CDA分为两个部分,即CDA第一步补助金和政府DollarforDollar 匹配。 根据CDA第一步补助金,父母将从政府现有捐款上限中初步获得3 000加元,存入孩子账户。) (0.7213630676269531)"
$ws.Range("E40").Value = "CCDA 第一步赠款是什么? (新加坡公民发展委员会第一阶段补助金是针对2016年3月24日出生的合格新加坡儿童的3 000项新补助金,支付给新加坡公民发展委员会子女,并构成政府现有总缴款上限的一部分。 新加坡公民在2016年3月24日出生的2016年3月24日出生的儿童,如果是新加坡公民,或在年满12岁之前成为公民,且有合法已婚父母,有资格享受新加坡公民发展委员会的福利。) (0.6798275709152222)"

# Row 41: "什麼是臍帶血幹細胞？"
$ws.Range("A41").Value = "什麼是臍帶血幹細胞？"
$ws.Range("C41").Value = "什么是绳索血干细胞? (一组干细胞在婴儿出生后,可以在细胞索中找到这些干细胞,具有独特的能力对血中发现的各种细胞类型进行区分,如红细胞下面的图表所描述的,红细胞下面的图表所传输的含氧白细胞可产生抗体和抗菌细菌,有助于血液凝聚。) (0.9091073870658875)"
$ws.Range("D41").Value = "绳索血干细胞能做什么? (细胞干细胞可替换和再生受损或患病的骨髓。 治疗与血液有关的癌症:固体肿瘤和免疫缺陷。 重新校正遗传缺陷:细胞疗法和再生医学。) (0.8465781807899475)"
$ws.Range("E41").Value = "使用绳索血干细胞可以治疗的疾病有哪些? (细胞干细胞目前正在作为一种标准治疗形式,治疗80多种疾病,如白血病淋巴瘤、地中海贫血和神经淋巴瘤,仅举几个例子。) (0.7302578687667847)"

# Row 42: "apa itu sel induk darah tali pusat?"
$ws.Range("A42").Value = "apa itu sel induk darah tali pusat?"
$ws.Range("C42").Value = "What are cord blood stem cells? (Cord blood stem cells are also known as Haematopoietic Stem Cells (HSCs), which are responsible for replenishing blood and regenerating the immune system. These stem cells can be found within the umbilical cord following the birth of a baby and have the unique ability to differentiate into various cell types found in blood as depicted in the diagram below:     
Red blood cells - which transport oxygen
White blood cells - which produce antibodies and fight bacteria
Platelets - which assist in blood clotting                                        ) (0.7981184720993042)"
$ws.Range("D42").Value = "What is cord blood? (Cord blood is blood that remains in the umbilical cord following the birth of a baby and after the umbilical cord is cut. During pregnancy, the umbilical cord functions as a lifeline between mother and child. After a baby�s delivery, the cord blood present in the umbilical cord could offer hope for the child as well as members of the family. Cord blood contains a rich source of lifesaving stem cells called Haematopoietic Stem Cells (HSCs).) (0.7619432806968689)"
$ws.Range("E42").Value = "What can cord blood stem cells do? (Cord blood stem cells can:
Replace and regenerate damaged or diseased bone marrow.
Treat blood related cancers, solid tumours and immunodeficiency disorders.
Correct genetic defects (sibling/allogeneic transplantation).
Potentially be used for cellular therapy and regenerative medicine.erative medicine.) (0.7380398511886597)"

# Row 43: "how do I enrol?"
$ws.Range("A43").Value = "how do I enrol?"
$ws.Range("C43").Value = "How do I enrol? (Contact us at our 24-Hour Hotline (65) 6238 0808 to make a non-obligatory appointment with our friendly Cord Blood Banking Consultants.                 *If you are due within the next 4 - 6 weeks, we strongly recommend you to contact us immediately for enrolment because this will ensure all necessary steps are completed before your baby's arrival.) (0.9978150725364685)"
$ws.Range("D43").Value = "What are the required documents to provide when I apply for my inter-bank GIRO Application? (You will be required to submit a copy of your child�s birth certificate, the Trustee�s name and your child�s Child Development Account (CDA) number to Cordlife.) (0.3560982644557953)"
$ws.Range("E43").Value = "What do I need to do if I wish to terminate my inter-bank GIRO used for my child�s Child Development Account (CDA) with Cordlife? (If you wish to terminate your inter-bank GIRO used for your child�s CDA with us, please contact the bank which you opened the CDA with, and inform Cordlife at 6238 0808 or email�customercare@cordlife.com.) (0.2969340682029724)"
